$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for changed rows ---

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.054.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.44%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.828.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.31%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.38%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6356"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.48%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.79"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.37%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2936"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.72%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07338"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.21%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.99%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07648"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.49%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.823.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.48%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.984"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.38%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6633"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.17%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.02%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008673"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.11%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.057"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.32%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "28.904.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.80%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.078.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.70%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.06%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.63%  "

# Row 23
$ws.Range("E23").Value = "  -0.10%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.108"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.17%  "

# Row 25
$ws.Range("E25").Value = "  -0.01%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.57%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.453"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.88%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1371"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.43%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.06%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.506"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.05%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.094"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.14%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.026"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.06%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.203"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.53%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05291"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.21%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.833"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.82%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7383"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.36%  "

# Row 37
$ws.Range("E37").Value = "  +2.21%  "

# Row 38
$ws.Range("E38").Value = "  -1.16%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.292.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.29%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.289"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.12%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8970"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.10%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.14%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.21%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.977.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.36%  "

# Row 47
$ws.Range("E47").Value = "  -0.47%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.13%  "

# Row 49
$ws.Range("E49").Value = "  -9.09%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.729"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.92%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05827"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.74%  "

# --- Rows 40 and 41 swap coin identity (VeChain <-> MXToken) with new data ---
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.749"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.21%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01781"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.77%  "
